# Fixed EROR Title is more than 31 characters
# Rename sheets so titles are <= 31 characters:
#   "<dia>-Asistente de Self Checkout" -> "<dia>-Self Checkout"
#   "<dia>-Representante de Servicio" -> "<dia>-RS"

$wb = $excel.ActiveWorkbook

$renameMap = @{
    "lunes-Asistente de Self Checkout"      = "lunes-Self Checkout"
    "lunes-Representante de Servicio"       = "lunes-RS"
    "martes-Asistente de Self Checkout"     = "martes-Self Checkout"
    "martes-Representante de Servicio"      = "martes-RS"
    "miércoles-Asistente de Self Checkout"  = "miércoles-Self Checkout"
    "miércoles-Representante de Servicio"   = "miércoles-RS"
    "jueves-Asistente de Self Checkout"     = "jueves-Self Checkout"
    "jueves-Representante de Servicio"      = "jueves-RS"
    "viernes-Asistente de Self Checkout"    = "viernes-Self Checkout"
    "viernes-Representante de Servicio"     = "viernes-RS"
    "sábado-Asistente de Self Checkout"     = "sábado-Self Checkout"
    "sábado-Representante de Servicio"      = "sábado-RS"
    "domingo-Asistente de Self Checkout"    = "domingo-Self Checkout"
    "domingo-Representante de Servicio"     = "domingo-RS"
}

foreach ($ws in $wb.Worksheets) {
    $oldName = $ws.Name
    if ($renameMap.ContainsKey($oldName)) {
        $ws.Name = $renameMap[$oldName]
    }
}
